$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Update the text of the run following " ... Jobbet med att separera denna
#    funktionalitet " (the run that starts with "var inte särskilt
#    krävande...") to the new, expanded wording about the synchronisation
#    unit.
# ---------------------------------------------------------------------------
$oldText1 = "var inte särskilt krävande. Dock behöves validering utav designen för att hantera packet som skapar upp fysikkroppar då nätverk och gameplay ligger på olika trådar och fysikkroppar kan inte skapas när gameplay uppdatera fysiken. Detta medför att saker måste läggas i temp arrayer vilket inte känns bra."
$newText1 = "var inte särskilt krävande. Dock behövdes en synkroniseringsenhet för att säkerhetsställa att paket som skickas via nätverket för att skapa upp fysikkroppar behandlas under uppdateringen av gameplaytråden. Detta p.g.a. att fysikkroppar inte får skapas upp under en eventuell uppdatering utav fysikvärlden, något vi inte kan säkerhetsställa på annat vis eftersom nätverk och gameplay ligger på separerade trådar."

$findRng1 = $d.Content
$found1 = $findRng1.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not locate the 'var inte sarskilt kravande...' run text"
}
$target1 = $d.Range($findRng1.Start, $findRng1.End)
$target1.Text = $newText1

# ---------------------------------------------------------------------------
# 2) The short "SyncedUpdate...." paragraph becomes a full sentence about the
#    refactoring outcome. The paragraph also keeps its "_GoBack" bookmark
#    (right where it was, immediately after the sentence) and gains a new
#    trailing run containing a single space placed after the bookmark.
# ---------------------------------------------------------------------------
$oldText2 = "SyncedUpdate...."
$newText2 = "Denna refaktorisering ledde till en mer komplett design i helthelt, men fler klasser var tvunget att läggas till vilket kan resultera i ett försvårande av användandet utav designen."

$findRng2 = $d.Content
$found2 = $findRng2.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not locate the 'SyncedUpdate....' paragraph text"
}
$target2Start = $findRng2.Start
$target2 = $d.Range($findRng2.Start, $findRng2.End)
# Write the new sentence plus the trailing space in one shot so the
# bookmark we add next lands inside existing run text (not at its very
# edge), which keeps it from being silently dropped.
$target2.Text = $newText2 + " "

$bmPos = $target2Start + $newText2.Length
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null

# ---------------------------------------------------------------------------
# 3) One of the three trailing empty paragraphs is removed (three become
#    two). Word will not allow deleting the very last paragraph mark in the
#    document, so remove the first empty paragraph right after the sentence
#    we just edited.
# ---------------------------------------------------------------------------
$targetParaIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Denna refaktorisering")) {
        $targetParaIndex = $i
        break
    }
}
if ($null -eq $targetParaIndex) {
    throw "Could not locate the updated 'Denna refaktorisering...' paragraph"
}
$d.Paragraphs($targetParaIndex + 1).Range.Delete()
